$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("F3").Value = 1.75
$ws.Range("J3").Value = 1.94
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AL3").Value = 50

# Row 4
$ws.Range("G4").Value = 1.5
$ws.Range("H4").Value = 6.8
$ws.Range("V4").Value = 1.13
$ws.Range("W4").Value = 2.96

# Row 5
$ws.Range("I5").Value = 2.46
$ws.Range("J5").Value = 1.2
$ws.Range("N5").Value = 1.11
$ws.Range("T5").Value = 1.03
$ws.Range("U5").Value = 1.03

# Row 6
$ws.Range("F6").Value = 1.12
$ws.Range("N6").Value = 1.32
$ws.Range("P6").Value = 1.32

# Row 7
$ws.Range("J7").Value = 1.2
$ws.Range("N7").Value = 1.11
$ws.Range("Q7").Value = 1.3
$ws.Range("R7").Value = 1.1
$ws.Range("S7").Value = 1.3
$ws.Range("T7").Value = 1.03
$ws.Range("U7").Value = 1.03

# Row 8
$ws.Range("J8").Value = 1.2
$ws.Range("N8").Value = 1.36
$ws.Range("P8").Value = 1.36

# Row 9
$ws.Range("G9").Value = 600
$ws.Range("J9").Value = 1.2
$ws.Range("N9").Value = 1.32
$ws.Range("P9").Value = 1.32

# Row 10
$ws.Range("J10").Value = 1.09
$ws.Range("N10").Value = 1.32
$ws.Range("P10").Value = 1.32
$ws.Range("R10").Value = 1.1
$ws.Range("T10").Value = 1.03
$ws.Range("U10").Value = 1.03

# Row 11
$ws.Range("H11").Value = 5.5
$ws.Range("J11").Value = 3.35
$ws.Range("L11").Value = 1.01
$ws.Range("M11").Value = 1.01
$ws.Range("N11").Value = 1.01
$ws.Range("O11").Value = 1.01
$ws.Range("R11").Value = 1.09
$ws.Range("S11").Value = 4.8
$ws.Range("T11").Value = 1.01
$ws.Range("U11").Value = 1.01
$ws.Range("V11").Value = 1.03
$ws.Range("W11").Value = 2.16
$ws.Range("X11").Value = 1000
$ws.Range("Y11").Value = 1000
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 1000
$ws.Range("AG11").Value = 1000
$ws.Range("AH11").Value = 1000
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("AO11").Value = 1000

# Row 12
$ws.Range("F12").Value = 2.78
$ws.Range("G12").Value = 3.15
$ws.Range("H12").Value = 2.8
$ws.Range("I12").Value = 3.15
$ws.Range("J12").Value = 2.9
$ws.Range("K12").Value = 3.3
$ws.Range("L12").Value = 1.01
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 2.62
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 1.54
$ws.Range("Q12").Value = 2.26
$ws.Range("R12").Value = 1.2
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = 2.02
$ws.Range("U12").Value = 1.81
$ws.Range("V12").Value = 1.46
$ws.Range("W12").Value = 1.46
$ws.Range("X12").Value = 980
$ws.Range("Y12").Value = 11
$ws.Range("Z12").Value = 22
$ws.Range("AA12").Value = 55
$ws.Range("AB12").Value = 11
$ws.Range("AC12").Value = 980
$ws.Range("AD12").Value = 16.5
$ws.Range("AE12").Value = 980
$ws.Range("AF12").Value = 23
$ws.Range("AG12").Value = 16.5
$ws.Range("AH12").Value = 980
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 55
$ws.Range("AK12").Value = 55
$ws.Range("AL12").Value = 80
$ws.Range("AM12").Value = 190
$ws.Range("AN12").Value = 55
$ws.Range("AO12").Value = 55
